$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A149").Value = "Golang Architect / Principal Backend Architect Only Local to GA"
$ws.Range("B149").Value = "https://www.dice.com/job-detail/09805e53-d7d8-4d5d-a4aa-ae268865ce37?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang"
$ws.Range("C149").Value = "Atlanta, Georgia"
$ws.Range("D149").Value = "Third Party"
$ws.Range("E149").Value = "Depends on Experience"
$ws.Range("F149").Value = "3BEES TECHNOLOGIES INC"

$ws.Range("A150").Value = "Mid Level Software Engineer - Python. W2. Hybrid Chicago"
$ws.Range("B150").Value = "https://www.dice.com/job-detail/feb492ca-33a8-467f-91d6-fa3f2c3765c0?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang"
$ws.Range("C150").Value = "Chicago, Illinois"
$ws.Range("D150").Value = "Contract"
$ws.Range("E150").Value = "Up to $75"
$ws.Range("F150").Value = "HSK Technologies, Inc."
